# "finished another fair comparison"
# Update the amp1..amp6 (columns B:G) sample data on the
# fair_comparison_stats_mix sheet, the "count" row, and the two summary
# rows (max / ratio-to-max) that depend on it. The mid1..mid6 columns
# (H:M) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2 : count
# ---------------------------------------------------------------------
$ws.Range("B2").Value = 895
$ws.Range("C2").Value = 920
$ws.Range("D2").Value = 823
$ws.Range("E2").Value = 806
$ws.Range("F2").Value = 402
$ws.Range("G2").Value = 385

# ---------------------------------------------------------------------
# Row 3 : mean
# ---------------------------------------------------------------------
$ws.Range("B3").Value = 1.44047709497207
$ws.Range("C3").Value = 1.6084152173913
$ws.Range("D3").Value = 1.7151907654920999
$ws.Range("E3").Value = 2.0797431761786598
$ws.Range("F3").Value = 2.1275547263681598
$ws.Range("G3").Value = 1.3647974025973999

# ---------------------------------------------------------------------
# Row 4 : std
# ---------------------------------------------------------------------
$ws.Range("B4").Value = 0.156994960236467
$ws.Range("C4").Value = 0.32530622617931898
$ws.Range("D4").Value = 0.39586748662978299
$ws.Range("E4").Value = 0.330158305408778
$ws.Range("F4").Value = 0.406280418914162
$ws.Range("G4").Value = 0.159139517187705

# ---------------------------------------------------------------------
# Row 5 : min
# ---------------------------------------------------------------------
$ws.Range("B5").Value = 1.0669999999999999
$ws.Range("C5").Value = 1.1060000000000001
$ws.Range("D5").Value = 1.274
$ws.Range("E5").Value = 1.486
$ws.Range("F5").Value = 1.571
$ws.Range("G5").Value = 1.085

# ---------------------------------------------------------------------
# Row 6 : 25%
# ---------------------------------------------------------------------
$ws.Range("B6").Value = 1.3445
$ws.Range("C6").Value = 1.4157500000000001
$ws.Range("D6").Value = 1.4910000000000001
$ws.Range("E6").Value = 1.9339999999999999
$ws.Range("F6").Value = 1.8022499999999999
$ws.Range("G6").Value = 1.21

# ---------------------------------------------------------------------
# Row 7 : 50%
# ---------------------------------------------------------------------
$ws.Range("B7").Value = 1.43
$ws.Range("C7").Value = 1.4815
$ws.Range("D7").Value = 1.544
$ws.Range("E7").Value = 1.984
$ws.Range("F7").Value = 1.8979999999999999
$ws.Range("G7").Value = 1.3879999999999999

# ---------------------------------------------------------------------
# Row 8 : 75%
# ---------------------------------------------------------------------
$ws.Range("B8").Value = 1.5449999999999999
$ws.Range("C8").Value = 1.5765
$ws.Range("D8").Value = 1.6160000000000001
$ws.Range("E8").Value = 2.0467499999999998
$ws.Range("F8").Value = 2.6287500000000001
$ws.Range("G8").Value = 1.5129999999999999

# ---------------------------------------------------------------------
# Row 9 : max (only B, C, E change -- D, F, G stay put)
# ---------------------------------------------------------------------
$ws.Range("B9").Value = 1.927
$ws.Range("C9").Value = 2.4260000000000002
$ws.Range("E9").Value = 2.9740000000000002

# ---------------------------------------------------------------------
# Row 11 : max-of-the-maxes -- amp columns (B:G) now look at row 3
# (the mean row) instead of row 8 (the 75% row); mid columns (H:M)
# are unchanged and still look at row 8.
# ---------------------------------------------------------------------
$ws.Range("B11").Formula = "=MAX(B3:G3)"
$ws.Range("C11").Formula = "=MAX(B3:G3)"
$ws.Range("D11").Formula = "=MAX(B3:G3)"
$ws.Range("E11").Formula = "=MAX(B3:G3)"
$ws.Range("F11").Formula = "=MAX(B3:G3)"
$ws.Range("G11").Formula = "=MAX(B3:G3)"

# ---------------------------------------------------------------------
# Row 12 : ratio vs the max -- amp columns (B:G) now divide row 11 by
# row 3 instead of row 8; mid columns (H:M) keep dividing by row 8.
# ---------------------------------------------------------------------
$ws.Range("B12").Formula = "=B11/B3"
$ws.Range("C12").Formula = "=C11/C3"
$ws.Range("D12").Formula = "=D11/D3"
$ws.Range("E12").Formula = "=E11/E3"
$ws.Range("F12").Formula = "=F11/F3"
$ws.Range("G12").Formula = "=G11/G3"

# ---------------------------------------------------------------------
# Row 12 label : "compensation" -> "scale"
# ---------------------------------------------------------------------
$ws.Range("A12").Value = "scale"

# ---------------------------------------------------------------------
# View bookkeeping (best effort; cosmetic only)
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D52").Select() | Out-Null
